# 2nd review of manual completed
# Append two new effort-log rows (56 and 57) to the "effort" sheet.
#
# Note: the two new shared strings must be created in the same order as in
# the target workbook (D57's text first, then D56's text), so that the
# resulting shared-string table indexes them as:
#   47 -> "Manual: 2nd review completed"
#   48 -> "Manual: 2nd review continued"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 57 first, so its string becomes shared-string index 47.
$ws.Range("A57").Value = 41255
$ws.Range("B57").Value = 1.75
$ws.Range("D57").Value = "Manual: 2nd review completed"

# Row 56 second, so its string becomes shared-string index 48.
$ws.Range("A56").Value = 41254
$ws.Range("B56").Value = 1
$ws.Range("D56").Value = "Manual: 2nd review continued"

# Match the new active cell / selection recorded in the workbook.
$ws.Range("D56").Select() | Out-Null
